# LicenseInfoResolver: Add original license source to ResolvedLicenseInfo.
# Updates the "ResolvedLicense(...)" strings on the Gradle project sheet so
# that they carry `originalExpressions` instead of `sources`, matching the
# richer ResolvedLicenseInfo model.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Update the declared-license detail cells with the new text ---------
$ws2.Range("C12").Value2 = "ResolvedLicense(license=EPL-1.0, originalDeclaredLicenses=[Eclipse Public License 1.0], originalExpressions={DECLARED=[EPL-1.0]}, locations=[])"
$ws2.Range("C13").Value2 = "ResolvedLicense(license=Apache-2.0, originalDeclaredLicenses=[Apache License, Version 2.0], originalExpressions={DECLARED=[Apache-2.0]}, locations=[])"
$ws2.Range("C14").Value2 = "ResolvedLicense(license=Apache-2.0, originalDeclaredLicenses=[Apache License, Version 2.0], originalExpressions={DECLARED=[Apache-2.0]}, locations=[])"
$ws2.Range("C15").Value2 = "ResolvedLicense(license=BSD-3-Clause, originalDeclaredLicenses=[New BSD License], originalExpressions={DECLARED=[BSD-3-Clause]}, locations=[])"

# --- 2. Give those cells their own highlighted style ------------------------
# Copy the formatting already used for "pulled in" detail cells (font +
# border) and then recolor the fill to the light-blue tone used for
# "Declared Licenses" detail rows elsewhere in the report.
$fmtSource = $ws2.Range("F11")
$detailCells = $ws2.Range("C12:C15")
$fmtSource.Copy()
$detailCells.PasteSpecial(-4122)
$detailCells.Interior.Color = 15128749

# --- 3. Update the frozen-pane scroll position / selection on this sheet ---
$ws2.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 11
$ws2.Range("C15").Select() | Out-Null

$excel.CutCopyMode = 0
